$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list (price + 1h volume change) per upstream diff.
# D-column price text is forced via a temporary Text number format so
# Excel keeps literal strings like "1.00" / "0.0870" instead of silently
# re-typing them as numbers; the style is reset to Normal afterwards so
# the cell keeps its original (default) formatting.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '91.300.04'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +2.51%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.115.40'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.47%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '218.37'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.53%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '621.58'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.32%  '

$ws.Range("E7").Value = '  +25.87%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.374'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.74%  '

$ws.Range("E9").Value = '  +0.09%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '3.112.88'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.33%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.731'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +22.75%  '

$ws.Range("E12").Value = '  +6.01%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000252'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +3.91%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '34.52'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +6.58%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.46'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +3.00%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '91.090.92'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.90%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.689.61'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.65%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.101.48'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.44%  '

$ws.Range("E19").Value = '  +12.21%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0000217'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.56%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.07'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +4.34%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '440.33'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +3.77%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.82'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +6.28%  '

$ws.Range("E24").Value = '  +4.93%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '6.19'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +9.53%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '88.46'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +7.07%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '12.28'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +2.57%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '3.281.37'
$ws.Range("D28").Style = "Normal"

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.00'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.10%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.167'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.95%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '9.11'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +12.01%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '523.90'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +2.13%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.895'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -16.89%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.74'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.75%  '

$ws.Range("E35").Value = '  +13.13%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '7.08'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +3.70%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '23.87'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +6.78%  '

$ws.Range("E38").Value = '  +2.76%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.86'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +3.07%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0870'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +25.52%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '22.28'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.01%  '

$ws.Range("E42").Value = '  +0.01%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.155'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +17.58%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.399'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +9.34%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.93'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +5.03%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '149.38'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +2.09%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '43.93'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.05%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.30'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +6.39%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '167.92'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +2.54%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '4.24'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +7.05%  '
